$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: Breadth First Search topic
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Breadth First search in Graph"
$ws.Range("H11").Value = "BFS"

# Row 12: Depth First Search topic
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Depth First Search in Graph"
$ws.Range("H12").Value = "DFS"

$ws.Range("H12").Select()
